$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 322974
$ws.Range("D2").Value = 411573047
$ws.Range("C4").Value = 319
$ws.Range("D4").Value = 456707
$ws.Range("C10").Value = 117562
$ws.Range("D10").Value = 172259949
$ws.Range("C12").Value = 59844
$ws.Range("D12").Value = 86375361
$ws.Range("C16").Value = 4012
$ws.Range("D16").Value = 5693492
$ws.Range("C20").Value = 6764
$ws.Range("D20").Value = 9440454
$ws.Range("C22").Value = 77998
$ws.Range("D22").Value = 97219735
$ws.Range("C27").Value = 289
$ws.Range("D27").Value = 414647
$ws.Range("C28").Value = 32613
$ws.Range("D28").Value = 47737942
$ws.Range("C30").Value = 11544
$ws.Range("D30").Value = 16607449
$ws.Range("C32").Value = 5
$ws.Range("D32").Value = 7500
$ws.Range("C33").Value = 1561
$ws.Range("D33").Value = 2192807
$ws.Range("C35").Value = 1852
$ws.Range("D35").Value = 2615179
$ws.Range("C36").Value = 97732
$ws.Range("D36").Value = 122958638
$ws.Range("C42").Value = 904
$ws.Range("D42").Value = 1330685
$ws.Range("C44").Value = 44526
$ws.Range("D44").Value = 65254435
$ws.Range("C46").Value = 9197
$ws.Range("D46").Value = 13196716
$ws.Range("C51").Value = 2340
$ws.Range("D51").Value = 3268018
$ws.Range("C52").Value = 69450
$ws.Range("D52").Value = 87108726
$ws.Range("C59").Value = 28316
$ws.Range("D59").Value = 41525953
$ws.Range("C62").Value = 11200
$ws.Range("D62").Value = 16196623
$ws.Range("C64").Value = 1363
$ws.Range("D64").Value = 1904737
$ws.Range("C68").Value = 1494
$ws.Range("D68").Value = 2093084
$ws.Range("C70").Value = 20618
$ws.Range("D70").Value = 27008207
$ws.Range("C74").Value = 7624
$ws.Range("D74").Value = 11164724
$ws.Range("C76").Value = 5156
$ws.Range("D76").Value = 7486610
$ws.Range("C79").Value = 141733
$ws.Range("D79").Value = 176687278
$ws.Range("C83").Value = 430
$ws.Range("D83").Value = 627824
$ws.Range("C85").Value = 63803
$ws.Range("D85").Value = 93512484
$ws.Range("C88").Value = 29907
$ws.Range("D88").Value = 43265941
$ws.Range("C91").Value = 2868
$ws.Range("D91").Value = 4056084
$ws.Range("C92").Value = 33483
$ws.Range("D92").Value = 45379530
$ws.Range("C96").Value = 8106
$ws.Range("D96").Value = 11918050
$ws.Range("C98").Value = 7458
$ws.Range("D98").Value = 10822367
$ws.Range("C102").Value = 10197
$ws.Range("D102").Value = 15394555
$ws.Range("C104").Value = 2537
$ws.Range("D104").Value = 4100088
$ws.Range("C106").Value = 3413
$ws.Range("D106").Value = 5504979
$ws.Range("C108").Value = 158
$ws.Range("D108").Value = 256445
$ws.Range("C109").Value = 196
$ws.Range("D109").Value = 298032
$ws.Range("C110").Value = 142418
$ws.Range("D110").Value = 176119579
$ws.Range("C112").Value = 73
$ws.Range("D112").Value = 104144
$ws.Range("C116").Value = 52971
$ws.Range("D116").Value = 77645580
$ws.Range("C118").Value = 27303
$ws.Range("D118").Value = 39560507
$ws.Range("C122").Value = 2294
$ws.Range("D122").Value = 3222668
$ws.Range("C124").Value = 518530
$ws.Range("D124").Value = 684872080
$ws.Range("C126").Value = 215
$ws.Range("D126").Value = 316509
$ws.Range("C129").Value = 1382
$ws.Range("D129").Value = 2048486
$ws.Range("C131").Value = 209555
$ws.Range("D131").Value = 308059771
$ws.Range("C134").Value = 184904
$ws.Range("D134").Value = 268890806
$ws.Range("C137").Value = 2856
$ws.Range("D137").Value = 4011969
$ws.Range("C139").Value = 6453
$ws.Range("D139").Value = 9117366
$ws.Range("C142").Value = 45019
$ws.Range("D142").Value = 60106966
$ws.Range("C148").Value = 14162
$ws.Range("D148").Value = 20764513
$ws.Range("C149").Value = 3794
$ws.Range("D149").Value = 5471002
$ws.Range("C154").Value = 391
$ws.Range("D154").Value = 552263
$ws.Range("C155").Value = 17748
$ws.Range("D155").Value = 23460996
$ws.Range("C159").Value = 7236
$ws.Range("D159").Value = 10526263
$ws.Range("C161").Value = 5054
$ws.Range("D161").Value = 7274629
$ws.Range("C166").Value = 18252
$ws.Range("D166").Value = 29420621
$ws.Range("C167").Value = 1996
$ws.Range("D167").Value = 3244449
$ws.Range("C168").Value = 269
$ws.Range("D168").Value = 430101
$ws.Range("C172").Value = 88331
$ws.Range("D172").Value = 110425175
$ws.Range("C179").Value = 34007
$ws.Range("D179").Value = 49872577
$ws.Range("C181").Value = 13117
$ws.Range("D181").Value = 18954462
$ws.Range("C185").Value = 1677
$ws.Range("D185").Value = 2357296
$ws.Range("C187").Value = 240028
$ws.Range("D187").Value = 298326830
$ws.Range("C195").Value = 86993
$ws.Range("D195").Value = 127522657
$ws.Range("C198").Value = 33246
$ws.Range("D198").Value = 47855799
$ws.Range("C201").Value = 5133
$ws.Range("D201").Value = 7308277
$ws.Range("C204").Value = 4954
$ws.Range("D204").Value = 6860923
$ws.Range("C207").Value = 265735
$ws.Range("D207").Value = 328837164
$ws.Range("C210").Value = 17
$ws.Range("D210").Value = 23935
$ws.Range("C216").Value = 95601
$ws.Range("D216").Value = 139859523
$ws.Range("C219").Value = 51798
$ws.Range("D219").Value = 74864968
$ws.Range("C222").Value = 4689
$ws.Range("D222").Value = 6584531
$ws.Range("C225").Value = 5866
$ws.Range("D225").Value = 8116077
$ws.Range("C228").Value = 107072
$ws.Range("D228").Value = 133876653
$ws.Range("C235").Value = 49752
$ws.Range("D235").Value = 72882937
$ws.Range("C237").Value = 12534
$ws.Range("D237").Value = 18024997
$ws.Range("C239").Value = 1898
$ws.Range("D239").Value = 2720882
$ws.Range("C241").Value = 2555
$ws.Range("D241").Value = 3574909
$ws.Range("C242").Value = 259630
$ws.Range("D242").Value = 327842062
$ws.Range("C248").Value = 835
$ws.Range("D248").Value = 1226563
$ws.Range("C250").Value = 96316
$ws.Range("D250").Value = 141130937
$ws.Range("C253").Value = 65600
$ws.Range("D253").Value = 95068640
$ws.Range("C255").Value = 2419
$ws.Range("D255").Value = 3411238
$ws.Range("C258").Value = 4671
$ws.Range("D258").Value = 6563526
